$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-12: step, CHR_intensity, pulseWidth, pulsePeriod, PulseNum, offTime, delayTime, iterationNum
$data = @(
    @(1,  5, 100, 200, 25, 0, 30, 1),
    @(2,  5, 100, 200, 25, 0, 60, 1),
    @(3,  5, 100, 200, 25, 0, 60, 1),
    @(4,  5, 100, 200, 25, 0, 60, 1),
    @(5,  5, 100, 200, 25, 0, 60, 1),
    @(6, 15, 100, 200, 25, 0, 60, 1),
    @(7, 15, 100, 200, 25, 0, 60, 1),
    @(8, 15, 100, 200, 25, 0, 60, 1),
    @(9, 15, 100, 200, 25, 0, 60, 1),
    @(10,15, 100, 200, 25, 0, 60, 1),
    @(11, 0, 100, 200, 25, 0, 60, 1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Formula = "=(D$row*E$row+F$row)*H$row+G$row*1000"
    $row = $row + 1
}

$ws.Range("K12").Select()
